$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF") matching the style of the existing
# header cells (bold, centered, thin border) by copying H1's formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-66
$values = @(
    @(2, 7, 7),
    @(3, 7, 7),
    @(4, 7, 8),
    @(5, 8, 8),
    @(6, 8, 8),
    @(7, 8, 8),
    @(8, 7, 7),
    @(9, 9, 9),
    @(10, 5, 6),
    @(11, 6, 7),
    @(12, 7, 7),
    @(13, 7, 7),
    @(14, 7, 7),
    @(15, 8, 8),
    @(16, 7, 7),
    @(17, 7, 7),
    @(18, 7, 7),
    @(19, 6, 6),
    @(20, 6, 7),
    @(21, 5, 6),
    @(22, 6, 6),
    @(23, 6, 8),
    @(24, 5, 6),
    @(25, 7, 7),
    @(26, 4, 6),
    @(27, 5, 6),
    @(28, 9, 9),
    @(29, 6, 7),
    @(30, 5, 7),
    @(31, 6, 8),
    @(32, 10, 10),
    @(33, 7, 8),
    @(34, 7, 7),
    @(35, 7, 7),
    @(36, 5, 5),
    @(37, 3, 5),
    @(38, 5, 5),
    @(39, 8, 8),
    @(40, 5, 6),
    @(41, 6, 6),
    @(42, 7, 7),
    @(43, 8, 9),
    @(44, 5, 6),
    @(45, 6, 7),
    @(46, 4, 5),
    @(47, 7, 7),
    @(48, 7, 8),
    @(49, 6, 7),
    @(50, 9, 9),
    @(51, 5, 6),
    @(52, 8, 8),
    @(53, 8, 8),
    @(54, 4, 5),
    @(55, 6, 7),
    @(56, 8, 8),
    @(57, 10, 12),
    @(58, 4, 5),
    @(59, 6, 6),
    @(60, 8, 8),
    @(61, 1, 3),
    @(62, 6, 6),
    @(63, 8, 8),
    @(64, 4, 4),
    @(65, 6, 6),
    @(66, 6, 6)
)

foreach ($row in $values) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}
